$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.27%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.84%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.128"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.84%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'-0.45%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'0.27%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.82%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'2.94%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9011"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.25%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'12.61%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1772"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.99%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'3.50%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.04188"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-5.60%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.1051"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.39%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001247"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.39%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005897"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.00%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.358"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.11%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.3296"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.92%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.543"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-7.03%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'0.85%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D22").Value = "'0.04096"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.74%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001223"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.07%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.003996"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.38%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001300"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'6.17%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D38").Value = "'0.02392"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'3.31%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05180"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.55%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007761"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-2.33%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1299"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.59%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.006956"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'7.03%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.001950"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.10%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008544"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.67%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3074"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.87%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006872"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'5.09%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.36%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.01083"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'220.24%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004202"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-40.20%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.36%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.36%"
$ws.Range("E51").Style = "Normal"
